$d = $word.ActiveDocument

$d.Content.Find.Execute("40×56=2240", $true, $true, $false, $false, $false, $true, 1, $false, "59×15=885", 2) | Out-Null
$d.Content.Find.Execute("55×14=770", $true, $true, $false, $false, $false, $true, 1, $false, "38×55=2090", 2) | Out-Null
$d.Content.Find.Execute("46×26=1196", $true, $true, $false, $false, $false, $true, 1, $false, "35×20=700", 2) | Out-Null
$d.Content.Find.Execute("78×42=3276", $true, $true, $false, $false, $false, $true, 1, $false, "57×81=4617", 2) | Out-Null
$d.Content.Find.Execute("11×45=495", $true, $true, $false, $false, $false, $true, 1, $false, "58×99=5742", 2) | Out-Null
$d.Content.Find.Execute("39×91=3549", $true, $true, $false, $false, $false, $true, 1, $false, "19×23=437", 2) | Out-Null
$d.Content.Find.Execute("78×15=1170", $true, $true, $false, $false, $false, $true, 1, $false, "50×80=4000", 2) | Out-Null
$d.Content.Find.Execute("36×72=2592", $true, $true, $false, $false, $false, $true, 1, $false, "58×37=2146", 2) | Out-Null
$d.Content.Find.Execute("45×38=1710", $true, $true, $false, $false, $false, $true, 1, $false, "18×34=612", 2) | Out-Null
$d.Content.Find.Execute("46×86=3956", $true, $true, $false, $false, $false, $true, 1, $false, "31×52=1612", 2) | Out-Null
$d.Content.Find.Execute("14×70=980", $true, $true, $false, $false, $false, $true, 1, $false, "97×62=6014", 2) | Out-Null
$d.Content.Find.Execute("17×70=1190", $true, $true, $false, $false, $false, $true, 1, $false, "64×51=3264", 2) | Out-Null
$d.Content.Find.Execute("22×84=1848", $true, $true, $false, $false, $false, $true, 1, $false, "97×92=8924", 2) | Out-Null
$d.Content.Find.Execute("46×93=4278", $true, $true, $false, $false, $false, $true, 1, $false, "63×47=2961", 2) | Out-Null
$d.Content.Find.Execute("83×59=4897", $true, $true, $false, $false, $false, $true, 1, $false, "94×49=4606", 2) | Out-Null
$d.Content.Find.Execute("56×13=728", $true, $true, $false, $false, $false, $true, 1, $false, "94×67=6298", 2) | Out-Null
$d.Content.Find.Execute("65×18=1170", $true, $true, $false, $false, $false, $true, 1, $false, "96×12=1152", 2) | Out-Null
$d.Content.Find.Execute("91×31=2821", $true, $true, $false, $false, $false, $true, 1, $false, "73×76=5548", 2) | Out-Null
$d.Content.Find.Execute("21×87=1827", $true, $true, $false, $false, $false, $true, 1, $false, "43×80=3440", 2) | Out-Null
$d.Content.Find.Execute("48×30=1440", $true, $true, $false, $false, $false, $true, 1, $false, "49×16=784", 2) | Out-Null
$d.Content.Find.Execute("90×13=1170", $true, $true, $false, $false, $false, $true, 1, $false, "93×15=1395", 2) | Out-Null
$d.Content.Find.Execute("12×40=480", $true, $true, $false, $false, $false, $true, 1, $false, "58×94=5452", 2) | Out-Null
$d.Content.Find.Execute("90×97=8730", $true, $true, $false, $false, $false, $true, 1, $false, "94×60=5640", 2) | Out-Null
$d.Content.Find.Execute("73×31=2263", $true, $true, $false, $false, $false, $true, 1, $false, "33×57=1881", 2) | Out-Null
$d.Content.Find.Execute("25×87=2175", $true, $true, $false, $false, $false, $true, 1, $false, "90×14=1260", 2) | Out-Null
